$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column ("Price") values are plain text in the source data (coinranking scrape),
# but several look like plain decimals (e.g. "0.9995") that Excel would otherwise
# auto-convert to a number on assignment. Force text via NumberFormat, assign, then
# clear the format back off so the cell keeps its original (no explicit style) look.

$d2 = $ws.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "29.527.22"
$d2.ClearFormats()
$ws.Range("E2").Value = "  -0.52%  "

$d3 = $ws.Range("D3")
$d3.NumberFormat = "@"
$d3.Value = "1.850.18"
$d3.ClearFormats()
$ws.Range("E3").Value = "  -0.26%  "

$d4 = $ws.Range("D4")
$d4.NumberFormat = "@"
$d4.Value = "0.9995"
$d4.ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "

$d5 = $ws.Range("D5")
$d5.NumberFormat = "@"
$d5.Value = "242.72"
$d5.ClearFormats()
$ws.Range("E5").Value = "  -0.70%  "

$d6 = $ws.Range("D6")
$d6.NumberFormat = "@"
$d6.Value = "0.6323"
$d6.ClearFormats()
$ws.Range("E6").Value = "  -1.21%  "

$ws.Range("E7").Value = "  +0.02%  "

$d8 = $ws.Range("D8")
$d8.NumberFormat = "@"
$d8.Value = "47.84"
$d8.ClearFormats()
$ws.Range("E8").Value = "  +1.24%  "

$d9 = $ws.Range("D9")
$d9.NumberFormat = "@"
$d9.Value = "0.07555"
$d9.ClearFormats()
$ws.Range("E9").Value = "  +0.96%  "

$d10 = $ws.Range("D10")
$d10.NumberFormat = "@"
$d10.Value = "0.2975"
$d10.ClearFormats()
$ws.Range("E10").Value = "  +0.36%  "

$d11 = $ws.Range("D11")
$d11.NumberFormat = "@"
$d11.Value = "24.29"
$d11.ClearFormats()
$ws.Range("E11").Value = "  -0.31%  "

$d12 = $ws.Range("D12")
$d12.NumberFormat = "@"
$d12.Value = "0.07677"
$d12.ClearFormats()
$ws.Range("E12").Value = "  +0.28%  "

$d13 = $ws.Range("D13")
$d13.NumberFormat = "@"
$d13.Value = "1.883.82"
$d13.ClearFormats()
$ws.Range("E13").Value = "  +1.48%  "

$d14 = $ws.Range("D14")
$d14.NumberFormat = "@"
$d14.Value = "5.023"
$d14.ClearFormats()
$ws.Range("E14").Value = "  -0.10%  "

$d15 = $ws.Range("D15")
$d15.NumberFormat = "@"
$d15.Value = "0.6851"
$d15.ClearFormats()
$ws.Range("E15").Value = "  -0.61%  "

$d16 = $ws.Range("D16")
$d16.NumberFormat = "@"
$d16.Value = "83.76"
$d16.ClearFormats()
$ws.Range("E16").Value = "  +0.07%  "

$d17 = $ws.Range("D17")
$d17.NumberFormat = "@"
$d17.Value = "0.000009852"
$d17.ClearFormats()
$ws.Range("E17").Value = "  +2.40%  "

$d18 = $ws.Range("D18")
$d18.NumberFormat = "@"
$d18.Value = "2.131.25"
$d18.ClearFormats()
$ws.Range("E18").Value = "  +1.05%  "

$d19 = $ws.Range("D19")
$d19.NumberFormat = "@"
$d19.Value = "6.203"
$d19.ClearFormats()
$ws.Range("E19").Value = "  +2.53%  "

$d20 = $ws.Range("D20")
$d20.NumberFormat = "@"
$d20.Value = "29.572.79"
$d20.ClearFormats()
$ws.Range("E20").Value = "  -0.45%  "

$d21 = $ws.Range("D21")
$d21.NumberFormat = "@"
$d21.Value = "234.10"
$d21.ClearFormats()
$ws.Range("E21").Value = "  -0.66%  "

$d22 = $ws.Range("D22")
$d22.NumberFormat = "@"
$d22.Value = "12.51"
$d22.ClearFormats()
$ws.Range("E22").Value = "  -1.01%  "

$d24 = $ws.Range("D24")
$d24.NumberFormat = "@"
$d24.Value = "7.612"
$d24.ClearFormats()
$ws.Range("E24").Value = "  +2.13%  "

$d25 = $ws.Range("D25")
$d25.NumberFormat = "@"
$d25.Value = "1.001"
$d25.ClearFormats()
$ws.Range("E25").Value = "  +0.01%  "

$d26 = $ws.Range("D26")
$d26.NumberFormat = "@"
$d26.Value = "155.72"
$d26.ClearFormats()
$ws.Range("E26").Value = "  -1.58%  "

$d27 = $ws.Range("D27")
$d27.NumberFormat = "@"
$d27.Value = "0.1386"
$d27.ClearFormats()
$ws.Range("E27").Value = "  -1.89%  "

$d28 = $ws.Range("D28")
$d28.NumberFormat = "@"
$d28.Value = "8.432"
$d28.ClearFormats()
$ws.Range("E28").Value = "  -0.95%  "

$ws.Range("E29").Value = "  -0.80%  "

$d30 = $ws.Range("D30")
$d30.NumberFormat = "@"
$d30.Value = "1.484"
$d30.ClearFormats()
$ws.Range("E30").Value = "  -0.83%  "

$d31 = $ws.Range("D31")
$d31.NumberFormat = "@"
$d31.Value = "0.05839"
$d31.ClearFormats()
$ws.Range("E31").Value = "  -6.24%  "

$d32 = $ws.Range("D32")
$d32.NumberFormat = "@"
$d32.Value = "1.279"
$d32.ClearFormats()
$ws.Range("E32").Value = "  +0.21%  "

$d33 = $ws.Range("D33")
$d33.NumberFormat = "@"
$d33.Value = "4.112"
$d33.ClearFormats()
$ws.Range("E33").Value = "  -0.96%  "

$d34 = $ws.Range("D34")
$d34.NumberFormat = "@"
$d34.Value = "4.048"
$d34.ClearFormats()
$ws.Range("E34").Value = "  -0.99%  "

$d35 = $ws.Range("D35")
$d35.NumberFormat = "@"
$d35.Value = "1.895"
$d35.ClearFormats()
$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("E36").Value = "  +0.06%  "

$d37 = $ws.Range("D37")
$d37.NumberFormat = "@"
$d37.Value = "0.7166"
$d37.ClearFormats()
$ws.Range("E37").Value = "  -1.41%  "

$ws.Range("E38").Value = "  -0.66%  "

$d39 = $ws.Range("D39")
$d39.NumberFormat = "@"
$d39.Value = "2.796"
$d39.ClearFormats()
$ws.Range("E39").Value = "  -1.20%  "

$d40 = $ws.Range("D40")
$d40.NumberFormat = "@"
$d40.Value = "1.236.06"
$d40.ClearFormats()
$ws.Range("E40").Value = "  +2.91%  "

$d41 = $ws.Range("D41")
$d41.NumberFormat = "@"
$d41.Value = "0.01773"
$d41.ClearFormats()
$ws.Range("E41").Value = "  -0.44%  "

$d42 = $ws.Range("D42")
$d42.NumberFormat = "@"
$d42.Value = "0.9152"
$d42.ClearFormats()
$ws.Range("E42").Value = "  -0.61%  "

$d43 = $ws.Range("D43")
$d43.NumberFormat = "@"
$d43.Value = "6.140"
$d43.ClearFormats()
$ws.Range("E43").Value = "  -0.05%  "

$d44 = $ws.Range("D44")
$d44.NumberFormat = "@"
$d44.Value = "2.042.90"
$d44.ClearFormats()
$ws.Range("E44").Value = "  +1.29%  "

$d45 = $ws.Range("D45")
$d45.NumberFormat = "@"
$d45.Value = "0.9997"
$d45.ClearFormats()
$ws.Range("E45").Value = "  -0.07%  "

$d46 = $ws.Range("D46")
$d46.NumberFormat = "@"
$d46.Value = "101.87"
$d46.ClearFormats()
$ws.Range("E46").Value = "  -0.22%  "

$d47 = $ws.Range("D47")
$d47.NumberFormat = "@"
$d47.Value = "67.50"
$d47.ClearFormats()
$ws.Range("E47").Value = "  +1.73%  "

$d48 = $ws.Range("D48")
$d48.NumberFormat = "@"
$d48.Value = "7.292"
$d48.ClearFormats()
$ws.Range("E48").Value = "  +9.33%  "

$d49 = $ws.Range("D49")
$d49.NumberFormat = "@"
$d49.Value = "9.173"
$d49.ClearFormats()
$ws.Range("E49").Value = "  -0.17%  "

$ws.Range("E50").Value = "  -1.21%  "

$d51 = $ws.Range("D51")
$d51.NumberFormat = "@"
$d51.Value = "0.4033"
$d51.ClearFormats()
$ws.Range("E51").Value = "  -0.43%  "
